$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the id cells (column A) for all data rows (2-128), leaving the
# header (A1) and the name column (B) untouched. Clear() (not
# ClearContents()) drops the cell record entirely rather than leaving an
# empty styled stub behind.
$ws.Range("A2:A128").Clear()
